# Update the answer cells in the division-practice table.
# Each old value is unique within the document, so a plain
# Find/Replace (MatchWholeWord + MatchCase, Wrap=wdFindContinue) for
# each pair is sufficient and avoids any table/cell indexing.

$d = $word.ActiveDocument

$replacements = @(
    @("325÷5=65, 0", "125÷7=17, 6"),
    @("789÷3=263, 0", "888÷7=126, 6"),
    @("952÷2=476, 0", "359÷9=39, 8"),
    @("532÷2=266, 0", "601÷7=85, 6"),
    @("642÷7=91, 5", "968÷4=242, 0"),
    @("133÷3=44, 1", "566÷5=113, 1"),
    @("116÷2=58, 0", "684÷5=136, 4"),
    @("655÷2=327, 1", "176÷9=19, 5"),
    @("187÷4=46, 3", "405÷3=135, 0"),
    @("834÷8=104, 2", "701÷9=77, 8"),
    @("943÷9=104, 7", "237÷2=118, 1"),
    @("189÷7=27, 0", "966÷3=322, 0"),
    @("884÷8=110, 4", "221÷5=44, 1"),
    @("973÷7=139, 0", "711÷3=237, 0"),
    @("186÷2=93, 0", "396÷2=198, 0"),
    @("662÷5=132, 2", "483÷4=120, 3"),
    @("236÷3=78, 2", "357÷4=89, 1"),
    @("231÷5=46, 1", "318÷4=79, 2"),
    @("605÷4=151, 1", "404÷5=80, 4"),
    @("385÷4=96, 1", "622÷5=124, 2"),
    @("802÷8=100, 2", "175÷8=21, 7"),
    @("145÷8=18, 1", "477÷2=238, 1"),
    @("304÷2=152, 0", "832÷2=416, 0"),
    @("493÷4=123, 1", "910÷4=227, 2"),
    @("376÷9=41, 7", "764÷4=191, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($replacements.Count) answer cells"
